# Update the OS-version entries in the "MACHINES" table on slide 3
# (Table 8, graphicFrame id=8) of the network schematic slide.
#
#   Capstone Server : Linux 4.6        -> Linux 4.15
#   Target 1        : Linux 3.2 - 4.9  -> Linux 3.16
#   Target 2        : Linux 3.2 - 4.9  -> Linux 3.16   (no trailing space)
#   ELK Server      : Linux 4.6        -> Linux 4.15

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table

# Each data row's single cell holds three paragraphs: name, IP, OS version.
# The OS-version text is always the 3rd paragraph.

$row3 = $tbl.Rows.Item(3).Cells.Item(1).Shape.TextFrame.TextRange.Paragraphs()
$row3.Item(3).Text = "Linux 4.15"

$row4 = $tbl.Rows.Item(4).Cells.Item(1).Shape.TextFrame.TextRange.Paragraphs()
$row4.Item(3).Text = "Linux 3.16 "

$row5 = $tbl.Rows.Item(5).Cells.Item(1).Shape.TextFrame.TextRange.Paragraphs()
$row5.Item(3).Text = "Linux 3.16"

$row6 = $tbl.Rows.Item(6).Cells.Item(1).Shape.TextFrame.TextRange.Paragraphs()
$row6.Item(3).Text = "Linux 4.15"
